# Update "想去人数" (interest count) figures on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 5697
$wsExhibit.Range("F8").Value = 387

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 5697
$wsAll.Range("F9").Value = 387
